$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 156
$ws.Range("I4").Value = 156
$ws.Range("K4").Value = 156
$ws.Range("M4").Value = -42

$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("M92").ClearContents()

$ws.Range("H98").Value = 752
$ws.Range("I98").Value = 752
$ws.Range("K98").Value = 752
$ws.Range("M98").Value = 746

$ws.Range("H106").Value = 1366.3334
$ws.Range("I106").Value = 1366.3334
$ws.Range("K106").Value = 1366.3334
$ws.Range("M106").Value = -735.3334

$ws.Range("H111").Value = 417.83334
$ws.Range("I111").Value = 401.4
$ws.Range("J111").Value = 500
$ws.Range("K111").Value = 1204.2
$ws.Range("L111").Value = 1500
$ws.Range("M111").Value = 1862.8
$ws.Range("N111").Value = -7634

$ws.Range("H122").Value = 752
$ws.Range("I122").Value = 752
$ws.Range("K122").Value = 2256
$ws.Range("M122").Value = 194

$ws.Range("H125").Value = 3251.1428
$ws.Range("I125").Value = 2291
$ws.Range("J125").Value = 4979.4
$ws.Range("K125").Value = 20619
$ws.Range("L125").Value = 44814.6
$ws.Range("M125").Value = -18159
$ws.Range("N125").Value = -49734.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 100
$ws.Range("I3").Value = 100
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 100
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 15
$ws.Range("N3").ClearContents()

$ws.Range("H4").Value = 450
$ws.Range("I4").Value = 450
$ws.Range("K4").Value = 450
$ws.Range("M4").Value = -334

$ws.Range("H74").Value = 1358.4
$ws.Range("I74").Value = 1349.5
$ws.Range("J74").Value = 1394
$ws.Range("K74").Value = 1349.5
$ws.Range("L74").Value = 1394
$ws.Range("M74").Value = -475.5
$ws.Range("N74").Value = -3142

$ws.Range("H77").Value = 1358.4
$ws.Range("I77").Value = 1349.5
$ws.Range("J77").Value = 1394
$ws.Range("K77").Value = 6747.5
$ws.Range("L77").Value = 6970
$ws.Range("M77").Value = -2379.5
$ws.Range("N77").Value = -15706

$ws.Range("H88").Value = 2229.3076
$ws.Range("J88").Value = 2435.6365
$ws.Range("L88").Value = 2435.6365
$ws.Range("N88").Value = -3247.6365

$ws.Range("H91").Value = 2229.3076
$ws.Range("J91").Value = 2435.6365
$ws.Range("L91").Value = 2435.6365
$ws.Range("N91").Value = -5243.636500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2982.9333
$ws.Range("I58").Value = 2768.889
$ws.Range("K58").Value = 2768.889
$ws.Range("M58").Value = -2565.889

$ws.Range("H99").Value = 2496.3333
$ws.Range("J99").Value = 2494.5
$ws.Range("L99").Value = 2494.5
$ws.Range("N99").Value = -5490.5

$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()

$ws.Range("H126").Value = 2496.3333
$ws.Range("J126").Value = 2494.5
$ws.Range("L126").Value = 7483.5
$ws.Range("N126").Value = -12423.5

$ws.Range("H136").Value = 2982.9333
$ws.Range("I136").Value = 2768.889
$ws.Range("K136").Value = 8306.667000000001
$ws.Range("M136").Value = -5756.667000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 209.375
$ws.Range("I2").Value = 147
$ws.Range("J2").Value = 230.16667
$ws.Range("K2").Value = 882
$ws.Range("L2").Value = 1381.00002
$ws.Range("M2").Value = -769
$ws.Range("N2").Value = -1607.00002

$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

$ws.Range("H107").Value = 685.0714
$ws.Range("J107").Value = 711.25
$ws.Range("L107").Value = 2133.75
$ws.Range("N107").Value = -5973.75

$ws.Range("H137").Value = 2495
$ws.Range("J137").Value = 2495
$ws.Range("L137").Value = 7485
$ws.Range("N137").Value = -17685

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 557.0625
$ws.Range("I97").Value = 587.1539
$ws.Range("J97").Value = 426.66666
$ws.Range("K97").Value = 587.1539
$ws.Range("L97").Value = 426.66666
$ws.Range("M97").Value = -91.15390000000002
$ws.Range("N97").Value = -1418.66666

$ws.Range("H102").Value = 2666.0833
$ws.Range("I102").Value = 1545
$ws.Range("K102").Value = 1545
$ws.Range("M102").Value = 77

$ws.Range("H122").Value = 4424.5
$ws.Range("I122").Value = 4231.6665
$ws.Range("J122").Value = 5003
$ws.Range("K122").Value = 12694.9995
$ws.Range("L122").Value = 15009
$ws.Range("M122").Value = -10244.9995
$ws.Range("N122").Value = -19909

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H125").Value = 30326
$ws.Range("J125").Value = 30326
$ws.Range("L125").Value = 30326
$ws.Range("N125").Value = -35246

$ws.Range("H126").Value = 9416.333000000001
$ws.Range("I126").Value = 9416.333000000001
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 28248.999
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -25778.999
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 18848.75
$ws.Range("I7").Value = 20470.889
$ws.Range("K7").Value = 20470.889
$ws.Range("M7").Value = -20358.889

$ws.Range("H22").Value = 1181.75
$ws.Range("J22").Value = 1471.2
$ws.Range("L22").Value = 1471.2
$ws.Range("N22").Value = -2061.2

$ws.Range("H27").Value = 1181.75
$ws.Range("J27").Value = 1471.2
$ws.Range("L27").Value = 1471.2
$ws.Range("N27").Value = -1685.2

$ws.Range("H46").Value = 3824.95
$ws.Range("I46").Value = 2500
$ws.Range("K46").Value = 2500
$ws.Range("M46").Value = -2312

$ws.Range("H68").Value = 1970.6666
$ws.Range("J68").Value = 1700
$ws.Range("L68").Value = 1700
$ws.Range("N68").Value = -3198

$ws.Range("H71").Value = 1970.6666
$ws.Range("J71").Value = 1700
$ws.Range("L71").Value = 8500
$ws.Range("N71").Value = -15988

$ws.Range("H126").Value = 18848.75
$ws.Range("I126").Value = 20470.889
$ws.Range("K126").Value = 61412.667
$ws.Range("M126").Value = -58942.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws.Range("H100").Value = 1542.4286
$ws.Range("I100").Value = 1699.5
$ws.Range("J100").Value = 600
$ws.Range("K100").Value = 3399
$ws.Range("L100").Value = 1200
$ws.Range("M100").Value = -2858
$ws.Range("N100").Value = -2282

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
